# "turn on the API" — add a new "period_value" column (all rows = 1) into
# the "SQL Query" sheet, inserted just before the existing "aggregation"
# column (which shifts from AG to AH, values unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AG, pushing the current AG ("aggregation") column
# and everything after it one position to the right (xlShiftToRight).
$ws.Columns("AG").Insert(-4161)

# Header for the newly inserted column.
$ws.Range("AG1").Value = "period_value"

# Populate the new column with 1 for every data row (2-157).
$ws.Range("AG2:AG157").Value = 1
